$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped timestamp on the previous (last) existing row.
$ws.Range("A17").Value2 = 45817.39392136574

# Append the newly scraped price data point as row 18.
$ws.Range("A18").Value2 = 45818.39371380122
$ws.Range("A18").NumberFormat = $ws.Range("A17").NumberFormat
$ws.Range("B18").Value = "EVOWHEY PROTEIN"
$ws.Range("C18").Value = "2Kg"
$ws.Range("D18").Value = "37,90€"
